$d = $word.ActiveDocument

# 1. Add the challenge name "ChatBOC" right after "Reto: "
$found1 = $d.Content.Find.Execute("Reto: ", $false, $false, $false, $false, $false, $true, 1, $false, "Reto: ChatBOC", 2)

# 2. Turn "... se constituye el EQUIPO 1:" into "... se constituye el: "Equipo A"."
$quoteOpen  = [char]0x201C
$quoteClose = [char]0x201D
$replacement2 = ": " + $quoteOpen + "Equipo A" + $quoteClose + "."
$found2 = $d.Content.Find.Execute(" EQUIPO 1:", $false, $false, $false, $false, $false, $true, 1, $false, $replacement2, 2)

# 3. Prefix the "MAIL" column header so it reads "EMAIL"
$found3 = $d.Content.Find.Execute("MAIL", $false, $false, $false, $false, $false, $true, 1, $false, "EMAIL", 2)

# 4. Remove the trailing blank row in the team-members table (after jblancog03@educantabria.es)
for ($i = 1; $i -le $d.Tables.Count; $i++) {
  $tbl = $d.Tables.Item($i)
  if ($tbl.Range.Text -like "*jblancog03@educantabria.es*") {
    for ($r = $tbl.Rows.Count; $r -ge 1; $r--) {
      $rowText = ($tbl.Rows.Item($r).Range.Text -replace "[\x07\x0d]", "")
      if ($rowText.Trim() -eq "") {
        $tbl.Rows.Item($r).Delete()
        break
      }
    }
    break
  }
}

Write-Output ("Reto found: " + $found1)
Write-Output ("Equipo found: " + $found2)
Write-Output ("Mail found: " + $found3)
